# Auto-generated script applying market-data refresh values
# per the commit 'chore: update Sheets via scheduled runner'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 8746.5
$ws.Range("I62").Value = 8746.5
$ws.Range("K62").Value = 8746.5
$ws.Range("M62").Value = -8122.5
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").Value = $null
# Row 65
$ws.Range("H65").Value = 8746.5
$ws.Range("I65").Value = 8746.5
$ws.Range("K65").Value = 43732.5
$ws.Range("M65").Value = -40612.5
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").Value = $null
# Row 98
$ws.Range("H98").Value = 3204.6316
$ws.Range("I98").Value = 3539.4707
$ws.Range("K98").Value = 3539.4707
$ws.Range("M98").Value = -2041.4707
# Row 103
$ws.Range("H103").Value = 2059.2
$ws.Range("I103").Value = 2074.25
$ws.Range("J103").Value = 1999
$ws.Range("K103").Value = 6222.75
$ws.Range("L103").Value = 5997
$ws.Range("M103").Value = -5636.75
$ws.Range("N103").Value = -7169
# Row 116
$ws.Range("H116").Value = 5450.4287
$ws.Range("I116").Value = 5450.4287
$ws.Range("K116").Value = 5450.4287
$ws.Range("M116").Value = -2008.4287
# Row 122
$ws.Range("H122").Value = 3204.6316
$ws.Range("I122").Value = 3539.4707
$ws.Range("K122").Value = 10618.4121
$ws.Range("M122").Value = -8168.4121
# Row 137
$ws.Range("H137").Value = 2638.4
$ws.Range("I137").Value = 1800
$ws.Range("J137").Value = 3197.3333
$ws.Range("K137").Value = 5400
$ws.Range("L137").Value = 9591.999899999999
$ws.Range("M137").Value = -2850
$ws.Range("N137").Value = -14691.9999
# Row 138
$ws.Range("H138").Value = 5552.0654
$ws.Range("J138").Value = 5577.1816
$ws.Range("L138").Value = 16731.5448
$ws.Range("N138").Value = -27011.5448

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4315.8335
$ws.Range("I61").Value = 4109.7
$ws.Range("K61").Value = 4109.7
$ws.Range("M61").Value = -3897.7
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null
# Row 111
$ws.Range("H111").Value = 644
$ws.Range("J111").Value = 644
$ws.Range("L111").Value = 644
$ws.Range("N111").Value = -8824
# Row 122
$ws.Range("H122").Value = 4499.5
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
# Row 132
$ws.Range("H132").Value = 1047.8334
$ws.Range("I132").Value = 857.6
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 2572.8
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -42.80000000000018
$ws.Range("N132").Value = -11057
# Row 136
$ws.Range("H136").Value = 4315.8335
$ws.Range("I136").Value = 4109.7
$ws.Range("K136").Value = 12329.1
$ws.Range("M136").Value = -9779.099999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 133
$ws.Range("H133").Value = 145000
$ws.Range("J133").Value = 145000
$ws.Range("L133").Value = 145000
$ws.Range("N133").Value = -155120

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3631.8
$ws.Range("I31").Value = 3232.8333
$ws.Range("J31").Value = 4230.25
$ws.Range("K31").Value = 3232.8333
$ws.Range("L31").Value = 4230.25
$ws.Range("M31").Value = -2937.8333
$ws.Range("N31").Value = -4820.25
# Row 34
$ws.Range("H34").Value = 3631.8
$ws.Range("I34").Value = 3232.8333
$ws.Range("J34").Value = 4230.25
$ws.Range("K34").Value = 3232.8333
$ws.Range("L34").Value = 4230.25
$ws.Range("M34").Value = -3030.8333
$ws.Range("N34").Value = -4634.25
# Row 62
$ws.Range("H62").Value = 12941.895
$ws.Range("I62").Value = 12110.883
$ws.Range("J62").Value = 20005.5
$ws.Range("K62").Value = 12110.883
$ws.Range("L62").Value = 20005.5
$ws.Range("M62").Value = -11486.883
$ws.Range("N62").Value = -21253.5
# Row 65
$ws.Range("H65").Value = 12941.895
$ws.Range("I65").Value = 12110.883
$ws.Range("J65").Value = 20005.5
$ws.Range("K65").Value = 60554.415
$ws.Range("L65").Value = 100027.5
$ws.Range("M65").Value = -57434.415
$ws.Range("N65").Value = -106267.5
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").Value = $null
# Row 132
$ws.Range("H132").Value = 5850.6665
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5850.6665
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17551.9995
$ws.Range("M132").Value = $null
$ws.Range("N132").Value = -22611.9995
# Row 134
$ws.Range("H134").Value = 4276.4287
$ws.Range("I134").Value = 3668.2144
$ws.Range("J134").Value = 5492.857
$ws.Range("K134").Value = 11004.6432
$ws.Range("L134").Value = 16478.571
$ws.Range("M134").Value = -8469.643199999999
$ws.Range("N134").Value = -21548.571

$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 1453.5
$ws.Range("I51").Value = 883.3333
$ws.Range("J51").Value = 1795.6
$ws.Range("K51").Value = 2649.9999
$ws.Range("L51").Value = 5386.799999999999
$ws.Range("M51").Value = -2189.9999
$ws.Range("N51").Value = -6306.799999999999
# Row 59
$ws.Range("H59").Value = 1200
$ws.Range("J59").Value = 1200
$ws.Range("L59").Value = 3600
$ws.Range("N59").Value = -4680
# Row 107
$ws.Range("H107").Value = 1272.5
$ws.Range("I107").Value = 1106.8182
$ws.Range("J107").Value = 1345.4
$ws.Range("K107").Value = 3320.4546
$ws.Range("L107").Value = 4036.2
$ws.Range("M107").Value = -1400.4546
$ws.Range("N107").Value = -7876.200000000001
# Row 131
$ws.Range("H131").Value = 20792.207
$ws.Range("I131").Value = 223699.6
$ws.Range("K131").Value = 671098.8
$ws.Range("M131").Value = -666058.8

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 6992.8945
$ws.Range("J80").Value = 8727.362999999999
$ws.Range("L80").Value = 8727.362999999999
$ws.Range("N80").Value = -10723.363
# Row 83
$ws.Range("H83").Value = 6992.8945
$ws.Range("J83").Value = 8727.362999999999
$ws.Range("L83").Value = 43636.815
$ws.Range("N83").Value = -53620.815
# Row 122
$ws.Range("H122").Value = 2569.1428
$ws.Range("I122").Value = 2584
$ws.Range("K122").Value = 7752
$ws.Range("M122").Value = -5302
# Row 126
$ws.Range("H126").Value = 5349.273
$ws.Range("J126").Value = 5967.8
$ws.Range("L126").Value = 17903.4
$ws.Range("N126").Value = -22843.4
# Row 132
$ws.Range("H132").Value = 3575.4375
$ws.Range("I132").Value = 2921.2
$ws.Range("K132").Value = 8763.599999999999
$ws.Range("M132").Value = -6233.599999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 6230.375
$ws.Range("I22").Value = 5731
$ws.Range("K22").Value = 5731
$ws.Range("M22").Value = -5436
# Row 27
$ws.Range("H27").Value = 6230.375
$ws.Range("I27").Value = 5731
$ws.Range("K27").Value = 5731
$ws.Range("M27").Value = -5624
# Row 68
$ws.Range("H68").Value = 2420.182
$ws.Range("I68").Value = 2477.5
$ws.Range("K68").Value = 2477.5
$ws.Range("M68").Value = -1728.5
# Row 71
$ws.Range("H71").Value = 2420.182
$ws.Range("I71").Value = 2477.5
$ws.Range("K71").Value = 12387.5
$ws.Range("M71").Value = -8643.5
# Row 86
$ws.Range("H86").Value = 85000
$ws.Range("J86").Value = 85000
$ws.Range("L86").Value = 85000
$ws.Range("N86").Value = -87372
# Row 89
$ws.Range("H89").Value = 85000
$ws.Range("J89").Value = 85000
$ws.Range("L89").Value = 255000
$ws.Range("N89").Value = -266856
# Row 116
$ws.Range("H116").Value = 99500
$ws.Range("J116").Value = 99500
$ws.Range("L116").Value = 99500
$ws.Range("N116").Value = -108678
# Row 132
$ws.Range("H132").Value = 3761.6667
$ws.Range("I132").Value = 2391.8333
$ws.Range("K132").Value = 7175.499899999999
$ws.Range("M132").Value = -4645.499899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 7148.4375
$ws.Range("I122").Value = 8589.833000000001
$ws.Range("J122").Value = 2824.25
$ws.Range("K122").Value = 25769.499
$ws.Range("L122").Value = 8472.75
$ws.Range("M122").Value = -23319.499
$ws.Range("N122").Value = -13372.75
# Row 132
$ws.Range("H132").Value = 11060.28
$ws.Range("I132").Value = 7341.6875
$ws.Range("J132").Value = 17671.111
$ws.Range("K132").Value = 22025.0625
$ws.Range("L132").Value = 53013.333
$ws.Range("M132").Value = -19495.0625
$ws.Range("N132").Value = -58073.333
